# Insert a new data row at row 117 (pushing existing rows 117..220 down to
# 118..221) and populate it with the new weekly price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(117).Insert()

$ws.Cells.Item(117, 1).Value = 10
$ws.Cells.Item(117, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(117, 3).Value = "La Araucanía"
$ws.Cells.Item(117, 4).Value = 44658
$ws.Cells.Item(117, 5).Value = 9
$ws.Cells.Item(117, 6).Value = "Fruta"
$ws.Cells.Item(117, 7).Value = 100102
$ws.Cells.Item(117, 8).Value = "Cítricos"
$ws.Cells.Item(117, 9).Value = 100102006
$ws.Cells.Item(117, 10).Value = "Pomelo"
$ws.Cells.Item(117, 11).Value = "Red Blush"
$ws.Cells.Item(117, 12).Value = "Primera"
$ws.Cells.Item(117, 13).Value = 40
$ws.Cells.Item(117, 14).Value = 15000
$ws.Cells.Item(117, 15).Value = 15000
$ws.Cells.Item(117, 16).Value = 15000
$ws.Cells.Item(117, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(117, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(117, 19).Value = 1000
$ws.Cells.Item(117, 20).Value = 15
